$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "future_climate"

$newSheet.Range("A1").Value = "ISEA3H09_WC30AS_V14_CMIP5_CCSM4_RCP85_2070_BIO.txt"
$newSheet.Range("A2").Value = "https://www.worldclim.org/data/bioclim.html"
$newSheet.Range("A3").Value = "WC30AS_V14: The 30 arc seconds climate surface data from Hijimans et al. 2005, available as the WorldClim database. This SDM is based on version 1.4"

$newSheet.Range("A4").Value = "CMIP5: The climate model being used"
$newSheet.Range("A5").Value = "CCSM4: Another climate model, Community Climate System Model 4"
$newSheet.Range("A6").Value = "RCP85: RCP8.5 - A global warming scenario that assumes nobody cooperated in preventing climate change"
$newSheet.Range("A7").Value = "2070: This is a prediction for 2070"

$newSheet.Range("A9").Value = "HID"
$newSheet.Range("B9").Value = "Hexagon ID"

$bioNames = @("BIO01_Mean","BIO02_Mean","BIO03_Mean","BIO04_Mean","BIO05_Mean","BIO06_Mean","BIO07_Mean","BIO08_Mean","BIO09_Mean","BIO10_Mean","BIO11_Mean","BIO12_Mean","BIO13_Mean","BIO14_Mean","BIO15_Mean","BIO16_Mean","BIO17_Mean","BIO18_Mean","BIO19_Mean")
$bioDescs = @(
  "Annual Mean Temperature",
  "Mean Diurnal Range (Mean of monthly (max temp - min temp))",
  "Isothermality (BIO2/BIO7) (×100)",
  "Temperature Seasonality (standard deviation ×100)",
  "Max Temperature of Warmest Month",
  "Min Temperature of Coldest Month",
  "Temperature Annual Range (BIO5-BIO6)",
  "Mean Temperature of Wettest Quarter",
  "Mean Temperature of Driest Quarter",
  "Mean Temperature of Warmest Quarter",
  "Mean Temperature of Coldest Quarter",
  "Annual Precipitation",
  "Precipitation of Wettest Month",
  "Precipitation of Driest Month",
  "Precipitation Seasonality (Coefficient of Variation)",
  "Precipitation of Wettest Quarter",
  "Precipitation of Driest Quarter",
  "Precipitation of Warmest Quarter",
  "Precipitation of Coldest Quarter"
)

for ($i = 0; $i -lt $bioNames.Length; $i++) {
    $row = 10 + $i
    $newSheet.Cells.Item($row, 1).Value = $bioNames[$i]
    $newSheet.Cells.Item($row, 2).Value = $bioDescs[$i]
}

$newSheet.Columns.Item(1).ColumnWidth = 14.285714285714286

$newSheet.Range("H23").Select() | Out-Null
